$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DOMA-1872: removed multi-tariff values for non-electricity meters
# from meter-import-example.xlsx
#
# Rows 2-5 (GVS/HVS) and rows 8-10 (TEPLO/GAZ) are single-tariff meters:
# set "Кол-во тарифов" (column F) to 1 and clear "Показание 2" / "Показание 3"
# (columns H and I). Rows 6-7 are electricity (EL) meters and keep their
# multi-tariff readings untouched.

$singleTariffRows = 2,3,4,5,8,9,10

foreach ($r in $singleTariffRows) {
    $ws.Range("F$r").Value = 1
    $ws.Range("H$r").ClearContents()
    $ws.Range("I$r").ClearContents()
}

# Also bump the gas meter number in row 10 from 22 to 33.
$ws.Range("E10").Value = 33
